$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook gains one new weekly sample (2 rows: "Primera" / "Segunda")
# inserted right before the current first Betarraga record (row 450),
# pushing all subsequent rows down by two. This also pushes the former
# last two rows (534/535) down to become the new last rows (536/537).
$ws.Rows.Item(450).Insert()
$ws.Rows.Item(450).Insert()

# Seed the two new rows with the same layout/formatting as the pair that
# was just pushed down to 452:453 (same Mercado/Region/Categoria/etc.),
# then overwrite the week-specific values (date, volume, prices).
$ws.Range("A452:R453").Copy()
$ws.Range("A450").PasteSpecial()

# Row 450 ("Primera")
$ws.Cells.Item(450, 4).Value2 = 45244
$ws.Cells.Item(450, 10).Value2 = 1600
$ws.Cells.Item(450, 11).Value2 = 550
$ws.Cells.Item(450, 12).Value2 = 600
$ws.Cells.Item(450, 13).Value2 = 575
$ws.Cells.Item(450, 16).Value2 = 192

# Row 451 ("Segunda")
$ws.Cells.Item(451, 4).Value2 = 45244
$ws.Cells.Item(451, 10).Value2 = 900
$ws.Cells.Item(451, 11).Value2 = 450
$ws.Cells.Item(451, 12).Value2 = 500
$ws.Cells.Item(451, 13).Value2 = 475
$ws.Cells.Item(451, 16).Value2 = 158
